$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6669418
$ws.Range("J17").Value = 7145669
$ws.Range("L17").Value = 21437007
$ws.Range("N17").Value = -21437343

$ws.Range("H86").Value = 7573.125
$ws.Range("I86").Value = 1258.6666
$ws.Range("K86").Value = 1258.6666
$ws.Range("M86").Value = -135.6666

$ws.Range("H89").Value = 7573.125
$ws.Range("I89").Value = 1258.6666
$ws.Range("K89").Value = 6293.333000000001
$ws.Range("M89").Value = -677.3330000000005

$ws.Range("H98").Value = 1153.6111
$ws.Range("I98").Value = 1020.38464
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 1020.38464
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 477.61536
$ws.Range("N98").Value = -4496

$ws.Range("H122").Value = 1153.6111
$ws.Range("I122").Value = 1020.38464
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3061.15392
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -611.1539199999997
$ws.Range("N122").Value = -9400

$ws.Range("H131").Value = 2191
$ws.Range("I131").Value = 1668.0625
$ws.Range("J131").Value = 2834.6155
$ws.Range("K131").Value = 5004.1875
$ws.Range("L131").Value = 8503.8465
$ws.Range("M131").Value = 35.8125
$ws.Range("N131").Value = -18583.8465

$ws.Range("H137").Value = 1930.8
$ws.Range("I137").Value = 1734.8182
$ws.Range("K137").Value = 5204.4546
$ws.Range("M137").Value = -2654.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16574
$ws.Range("I32").Value = 13815.077
$ws.Range("K32").Value = 13815.077
$ws.Range("M32").Value = -13528.077

$ws.Range("H74").Value = 29413484
$ws.Range("I74").Value = 43478876
$ws.Range("J74").Value = 4031.6365
$ws.Range("K74").Value = 43478876
$ws.Range("L74").Value = 4031.6365
$ws.Range("M74").Value = -43478002
$ws.Range("N74").Value = -5779.636500000001

$ws.Range("H77").Value = 29413484
$ws.Range("I77").Value = 43478876
$ws.Range("J77").Value = 4031.6365
$ws.Range("K77").Value = 217394380
$ws.Range("L77").Value = 20158.1825
$ws.Range("M77").Value = -217390012
$ws.Range("N77").Value = -28894.1825

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7399.8
$ws.Range("I20").Value = 7399.8
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 7399.8
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -7152.8
$ws.Range("N20").ClearContents()

$ws.Range("H47").Value = 119850
$ws.Range("J47").Value = 119850
$ws.Range("L47").Value = 119850
$ws.Range("N47").Value = -120890

$ws.Range("H86").Value = 1940.6666
$ws.Range("I86").Value = 1682.6471
$ws.Range("J86").Value = 2278.077
$ws.Range("K86").Value = 1682.6471
$ws.Range("L86").Value = 2278.077
$ws.Range("M86").Value = -559.6470999999999
$ws.Range("N86").Value = -4524.077

$ws.Range("H89").Value = 1940.6666
$ws.Range("I89").Value = 1682.6471
$ws.Range("J89").Value = 2278.077
$ws.Range("K89").Value = 8413.235499999999
$ws.Range("L89").Value = 11390.385
$ws.Range("M89").Value = -2797.235499999999
$ws.Range("N89").Value = -22622.385

$ws.Range("H134").Value = 3722.6765
$ws.Range("I134").Value = 3899.0938
$ws.Range("K134").Value = 11697.2814
$ws.Range("M134").Value = -9162.2814

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3035.18
$ws.Range("I31").Value = 1480.5883
$ws.Range("J31").Value = 6338.6875
$ws.Range("K31").Value = 1480.5883
$ws.Range("L31").Value = 6338.6875
$ws.Range("M31").Value = -1185.5883
$ws.Range("N31").Value = -6928.6875

$ws.Range("H34").Value = 3035.18
$ws.Range("I34").Value = 1480.5883
$ws.Range("J34").Value = 6338.6875
$ws.Range("K34").Value = 1480.5883
$ws.Range("L34").Value = 6338.6875
$ws.Range("M34").Value = -1278.5883
$ws.Range("N34").Value = -6742.6875

$ws.Range("H109").Value = 247517500
$ws.Range("J109").Value = 247517500
$ws.Range("L109").Value = 247517500
$ws.Range("N109").Value = -247519580

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1264.8
$ws.Range("J5").Value = 1918.7693
$ws.Range("L5").Value = 5756.3079
$ws.Range("N5").Value = -5980.3079

$ws.Range("H131").Value = 692.0404
$ws.Range("J131").Value = 718.6889
$ws.Range("L131").Value = 2156.0667
$ws.Range("N131").Value = -12236.0667

$ws.Range("H135").Value = 1264.8
$ws.Range("J135").Value = 1918.7693
$ws.Range("L135").Value = 17268.9237
$ws.Range("N135").Value = -22338.9237

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 6593333.5
$ws.Range("I12").Value = 6593333.5
$ws.Range("K12").Value = 6593333.5
$ws.Range("M12").Value = -6593193.5

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws.Range("H24").Value = 2859.8
$ws.Range("J24").Value = 2859.8
$ws.Range("L24").Value = 2859.8
$ws.Range("N24").Value = -3545.8

$ws.Range("H25").Value = 4336
$ws.Range("J25").Value = 4336
$ws.Range("L25").Value = 4336
$ws.Range("N25").Value = -4796

$ws.Range("H82").Value = 1611
$ws.Range("I82").Value = 1450
$ws.Range("K82").Value = 1450
$ws.Range("M82").Value = -1089

$ws.Range("H85").Value = 1611
$ws.Range("I85").Value = 1450
$ws.Range("K85").Value = 1450
$ws.Range("M85").Value = -202

$ws.Range("H93").Value = 2497.5
$ws.Range("I93").Value = 2995
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 2995
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -1747
$ws.Range("N93").Value = -4496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 30380
$ws.Range("J109").Value = 30380
$ws.Range("L109").Value = 30380
$ws.Range("N109").Value = -33154

$ws.Range("H115").Value = 28484.857
$ws.Range("I115").Value = 20000
$ws.Range("J115").Value = 29899
$ws.Range("K115").Value = 20000
$ws.Range("L115").Value = 29899
$ws.Range("M115").Value = -18433
$ws.Range("N115").Value = -33033

$ws.Range("H122").Value = 873.1177
$ws.Range("I122").Value = 608.9091
$ws.Range("J122").Value = 1357.5
$ws.Range("K122").Value = 1826.7273
$ws.Range("L122").Value = 4072.5
$ws.Range("M122").Value = 623.2727
$ws.Range("N122").Value = -8972.5
